$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add total scip time to the output log.
# Columns P/Q already held solution/runtime info for some rows; we now also
# add a new column R with a remark, and fill in P/Q for rows 5 and 16 plus
# the new remark for row 4.

# Row 5 (E-n51-k5 area): fill Q5, P5, then R5
$ws.Range("Q5").Value2 = "3300s"
$ws.Range("P5").Value2 = "803,3.."
$ws.Range("R5").Value2 = "1 non elementary path"

# Row 16: fill Q16, P16
$ws.Range("Q16").Value2 = "2020s"
$ws.Range("P16").Value2 = "14735,9…"

# Row 16: R16 reuses the same remark as R5
$ws.Range("R16").Value2 = "1 non elementary path"

# Row 4: P4/Q4 already had data values elsewhere in the workbook; set them
# here too, then add the new remark in R4 last.
$ws.Range("P4").Value2 = "484,0…"
$ws.Range("Q4").Value2 = "20000s"
$ws.Range("R4").Value2 = "elementary, sehr ungleiche Laufzeit"

# Restore the active selection to match the saved workbook state.
$ws.Range("F10").Select()
